# Auto-generated: applies 2025-12-14 violent crime data update
# Updates column L (year 2025) values across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets, per the commit diff.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 6317   # was 6291
$ws.Range('L3').Value = 6826   # was 6810
$ws.Range('L4').Value = 1699   # was 1694
$ws.Range('L5').Value = 403   # was 402
$ws.Range('L6').Value = 5615   # was 5598
$ws.Range('L7').Value = 20860   # was 20795

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 187   # was 186
$ws.Range('L6').Value = 168   # was 166
$ws.Range('L7').Value = 661   # was 659
$ws.Range('L8').Value = 1379   # was 1376
$ws.Range('L10').Value = 139   # was 136
$ws.Range('L11').Value = 346   # was 344
$ws.Range('L15').Value = 173   # was 170
$ws.Range('L18').Value = 144   # was 143
$ws.Range('L19').Value = 566   # was 565
$ws.Range('L20').Value = 528   # was 527
$ws.Range('L21').Value = 68   # was 67
$ws.Range('L22').Value = 70   # was 69
$ws.Range('L29').Value = 1160   # was 1157
$ws.Range('L31').Value = 209   # was 208
$ws.Range('L33').Value = 937   # was 935
$ws.Range('L36').Value = 267   # was 264
$ws.Range('L37').Value = 802   # was 798
$ws.Range('L38').Value = 22   # was 20
$ws.Range('L41').Value = 90   # was 89
$ws.Range('L42').Value = 663   # was 661
$ws.Range('L44').Value = 141   # was 140
$ws.Range('L46').Value = 47   # was 46
$ws.Range('L52').Value = 444   # was 442
$ws.Range('L54').Value = 447   # was 445
$ws.Range('L63').Value = 61   # was 68
$ws.Range('L64').Value = 130   # was 129
$ws.Range('L65').Value = 409   # was 407
$ws.Range('L67').Value = 726   # was 723
$ws.Range('L73').Value = 163   # was 162
$ws.Range('L76').Value = 327   # was 323
$ws.Range('L79').Value = 575   # was 571
$ws.Range('L85').Value = 1037   # was 1033
$ws.Range('L87').Value = 57   # was 56
$ws.Range('L88').Value = 221   # was 220
$ws.Range('L90').Value = 220   # was 219
$ws.Range('L91').Value = 280   # was 279
$ws.Range('L95').Value = 294   # was 293
$ws.Range('L96').Value = 232   # was 231
$ws.Range('L98').Value = 111   # was 109
$ws.Range('L101').Value = 20860   # was 20795

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L2').Value = 73   # was 72
$ws.Range('L7').Value = 232   # was 231

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 231   # was 230
$ws.Range('L3').Value = 209   # was 208
$ws.Range('L7').Value = 661   # was 659

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 129   # was 128
$ws.Range('L6').Value = 89   # was 88
$ws.Range('L7').Value = 346   # was 344

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 312   # was 311
$ws.Range('L3').Value = 432   # was 430
$ws.Range('L4').Value = 59   # was 58
$ws.Range('L7').Value = 1037   # was 1033

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 140   # was 138
$ws.Range('L7').Value = 444   # was 442

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L3').Value = 487   # was 486
$ws.Range('L6').Value = 334   # was 332
$ws.Range('L7').Value = 1379   # was 1376

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 255   # was 254
$ws.Range('L3').Value = 331   # was 330
$ws.Range('L7').Value = 937   # was 935

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L3').Value = 97   # was 96
$ws.Range('L7').Value = 294   # was 293

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 243   # was 241
$ws.Range('L3').Value = 285   # was 284
$ws.Range('L6').Value = 207   # was 206
$ws.Range('L7').Value = 802   # was 798

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 149   # was 147
$ws.Range('L7').Value = 409   # was 407

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('L6').Value = 54   # was 53
$ws.Range('L7').Value = 209   # was 208

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 207   # was 205
$ws.Range('L6').Value = 169   # was 168
$ws.Range('L7').Value = 726   # was 723

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L4').Value = 37   # was 36
$ws.Range('L6').Value = 214   # was 213
$ws.Range('L7').Value = 447   # was 445

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 349   # was 348
$ws.Range('L3').Value = 447   # was 445
$ws.Range('L7').Value = 1160   # was 1157

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L4').Value = 28   # was 27
$ws.Range('L7').Value = 566   # was 565

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('L2').Value = 55   # was 54
$ws.Range('L7').Value = 141   # was 140

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L4').Value = 41   # was 39
$ws.Range('L6').Value = 146   # was 144
$ws.Range('L7').Value = 327   # was 323

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L2').Value = 66   # was 64
$ws.Range('L7').Value = 168   # was 166

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L2').Value = 27   # was 26
$ws.Range('L7').Value = 90   # was 89

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('L2').Value = 180   # was 179
$ws.Range('L5').Value = 17   # was 16
$ws.Range('L7').Value = 663   # was 661

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L2').Value = 54   # was 53
$ws.Range('L6').Value = 39   # was 37
$ws.Range('L7').Value = 139   # was 136

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('L6').Value = 16   # was 15
$ws.Range('L7').Value = 47   # was 46

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L2').Value = 97   # was 96
$ws.Range('L7').Value = 280   # was 279

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('L3').Value = 18   # was 17
$ws.Range('L7').Value = 68   # was 67

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L2').Value = 179   # was 178
$ws.Range('L3').Value = 184   # was 182
$ws.Range('L6').Value = 153   # was 152
$ws.Range('L7').Value = 575   # was 571

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L3').Value = 39   # was 38
$ws.Range('L7').Value = 130   # was 129

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 166   # was 165
$ws.Range('L7').Value = 528   # was 527

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('L2').Value = 51   # was 50
$ws.Range('L7').Value = 144   # was 143

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L2').Value = 91   # was 90
$ws.Range('L3').Value = 87   # was 86
$ws.Range('L4').Value = 21   # was 20
$ws.Range('L7').Value = 267   # was 264

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('L3').Value = 56   # was 55
$ws.Range('L4').Value = 15   # was 14
$ws.Range('L6').Value = 35   # was 34
$ws.Range('L7').Value = 173   # was 170

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('L4').Value = 14   # was 13
$ws.Range('L6').Value = 50   # was 49
$ws.Range('L7').Value = 111   # was 109

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L3').Value = 50   # was 49
$ws.Range('L7').Value = 163   # was 162

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L6').Value = 49   # was 48
$ws.Range('L7').Value = 187   # was 186

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('L4').Value = 14   # was 13
$ws.Range('L7').Value = 221   # was 220

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('L6').Value = 63   # was 62
$ws.Range('L7').Value = 220   # was 219

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range('L2').Value = 23   # was 22
$ws.Range('L7').Value = 70   # was 69

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('L6').Value = 21   # was 20
$ws.Range('L7').Value = 57   # was 56

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range('L3').Value = 7   # was 5
$ws.Range('L6').Value = 22   # was 20
